# Auto-generated Excel COM-interop script
# Applies the Titan_Profits market-data refresh diff (scheduled runner update)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1884.5625
$ws.Range("I40").Value = 1910.2
$ws.Range("J40").Value = 1872.909
$ws.Range("K40").Value = 1910.2
$ws.Range("L40").Value = 1872.909
$ws.Range("M40").Value = -1735.2
$ws.Range("N40").Value = -2222.909
# Row 51
$ws.Range("H51").Value = 3465.9722
$ws.Range("I51").Value = 2800.2
$ws.Range("J51").Value = 3573.3547
$ws.Range("K51").Value = 2800.2
$ws.Range("L51").Value = 3573.3547
$ws.Range("M51").Value = -2316.2
$ws.Range("N51").Value = -4541.3547
# Row 76
$ws.Range("H76").Value = 3971175.8
$ws.Range("I76").Value = 4447318.5
$ws.Range("J76").Value = 3320
$ws.Range("K76").Value = 4447318.5
$ws.Range("L76").Value = 3320
$ws.Range("M76").Value = -4447003.5
$ws.Range("N76").Value = -3950
# Row 79
$ws.Range("H79").Value = 3971175.8
$ws.Range("I79").Value = 4447318.5
$ws.Range("J79").Value = 3320
$ws.Range("K79").Value = 4447318.5
$ws.Range("L79").Value = 3320
$ws.Range("M79").Value = -4446226.5
$ws.Range("N79").Value = -5504
# Row 112
$ws.Range("H112").Value = 11858866
$ws.Range("J112").Value = 12988187
$ws.Range("L112").Value = 38964561
$ws.Range("N112").Value = -38966777
# Row 121
$ws.Range("H121").Value = 668.5833
$ws.Range("J121").Value = 676.25714
$ws.Range("L121").Value = 2028.77142
$ws.Range("N121").Value = -5522.77142
# Row 138
$ws.Range("H138").Value = 9530780
$ws.Range("I138").Value = 5498657
$ws.Range("J138").Value = 10646049
$ws.Range("K138").Value = 16495971
$ws.Range("L138").Value = 31938147
$ws.Range("M138").Value = -16490831
$ws.Range("N138").Value = -31948427

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 59
$ws.Range("H59").Value = 7400
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
# Row 61
$ws.Range("H61").Value = 2280.1333
$ws.Range("I61").Value = 1871.2307
$ws.Range("J61").Value = 4938
$ws.Range("K61").Value = 1871.2307
$ws.Range("L61").Value = 4938
$ws.Range("M61").Value = -1659.2307
$ws.Range("N61").Value = -5362
# Row 136
$ws.Range("H136").Value = 2280.1333
$ws.Range("I136").Value = 1871.2307
$ws.Range("J136").Value = 4938
$ws.Range("K136").Value = 5613.6921
$ws.Range("L136").Value = 14814
$ws.Range("M136").Value = -3063.6921
$ws.Range("N136").Value = -19914

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2456.4285
$ws.Range("I86").Value = 2532.5
$ws.Range("K86").Value = 2532.5
$ws.Range("M86").Value = -1409.5
# Row 89
$ws.Range("H89").Value = 2456.4285
$ws.Range("I89").Value = 2532.5
$ws.Range("K89").Value = 12662.5
$ws.Range("M89").Value = -7046.5
# Row 105
$ws.Range("H105").Value = 3209.1667
$ws.Range("I105").Value = 3033.2
$ws.Range("J105").Value = 3609.0908
$ws.Range("K105").Value = 3033.2
$ws.Range("L105").Value = 3609.0908
$ws.Range("M105").Value = -1286.2
$ws.Range("N105").Value = -7103.0908
# Row 134
$ws.Range("H134").Value = 7920568
$ws.Range("I134").Value = 1841207.2
$ws.Range("J134").Value = 20839208
$ws.Range("K134").Value = 5523621.6
$ws.Range("L134").Value = 62517624
$ws.Range("M134").Value = -5521086.6
$ws.Range("N134").Value = -62522694

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 87.26087
$ws.Range("I12").Value = 113.09091
$ws.Range("J12").Value = 63.583332
$ws.Range("K12").Value = 339.27273
$ws.Range("L12").Value = 190.749996
$ws.Range("M12").Value = -166.27273
$ws.Range("N12").Value = -536.749996
# Row 13
$ws.Range("H13").Value = 331.22223
$ws.Range("I13").Value = 96.833336
$ws.Range("K13").Value = 290.500008
$ws.Range("M13").Value = -122.500008
# Row 25
$ws.Range("H25").Value = 866
$ws.Range("I25").Value = 99
$ws.Range("J25").Value = 1249.5
$ws.Range("K25").Value = 297
$ws.Range("L25").Value = 3748.5
$ws.Range("M25").Value = -128
$ws.Range("N25").Value = -4086.5
# Row 30
$ws.Range("H30").Value = 866
$ws.Range("I30").Value = 99
$ws.Range("J30").Value = 1249.5
$ws.Range("K30").Value = 297
$ws.Range("L30").Value = 3748.5
$ws.Range("M30").Value = -195
$ws.Range("N30").Value = -3952.5
# Row 80
$ws.Range("H80").Value = 1093.2
# Row 83
$ws.Range("H83").Value = 1093.2
# Row 114
$ws.Range("H114").Value = 437.73685
$ws.Range("I114").Value = 150.45454
$ws.Range("J114").Value = 832.75
$ws.Range("K114").Value = 451.36362
$ws.Range("L114").Value = 2498.25
$ws.Range("M114").Value = 2802.63638
$ws.Range("N114").Value = -9006.25
# Row 117
$ws.Range("H117").Value = 661.1429000000001
$ws.Range("I117").Value = 282
$ws.Range("J117").Value = 1166.6666
$ws.Range("K117").Value = 846
$ws.Range("L117").Value = 3499.9998
$ws.Range("M117").Value = 2596
$ws.Range("N117").Value = -10383.9998
# Row 129
$ws.Range("H129").Value = 1133.95
$ws.Range("I129").Value = 414.77777
$ws.Range("J129").Value = 1722.3636
$ws.Range("K129").Value = 1244.33331
$ws.Range("L129").Value = 5167.0908
$ws.Range("M129").Value = 3755.66669
$ws.Range("N129").Value = -15167.0908
# Row 132
$ws.Range("H132").Value = 1555.5555
$ws.Range("I132").Value = 1700
$ws.Range("J132").Value = 1537.5
$ws.Range("K132").Value = 15300
$ws.Range("L132").Value = 13837.5
$ws.Range("M132").Value = -12770
$ws.Range("N132").Value = -18897.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
# Row 70
$ws.Range("H70").Value = 6199.4
$ws.Range("I70").Value = 7037.5386
$ws.Range("J70").Value = 4642.857
$ws.Range("K70").Value = 7037.5386
$ws.Range("L70").Value = 4642.857
$ws.Range("M70").Value = -6767.5386
$ws.Range("N70").Value = -5182.857
# Row 73
$ws.Range("H73").Value = 6199.4
$ws.Range("I73").Value = 7037.5386
$ws.Range("J73").Value = 4642.857
$ws.Range("K73").Value = 7037.5386
$ws.Range("L73").Value = 4642.857
$ws.Range("M73").Value = -6101.5386
$ws.Range("N73").Value = -6514.857
# Row 132
$ws.Range("H132").Value = 3759.9546
$ws.Range("I132").Value = 3567.0557
$ws.Range("K132").Value = 10701.1671
$ws.Range("M132").Value = -8171.167099999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 13
$ws.Range("H13").Value = 31000
$ws.Range("I13").Value = 44000
$ws.Range("J13").Value = 5000
$ws.Range("K13").Value = 44000
$ws.Range("L13").Value = 5000
$ws.Range("M13").Value = -43860
$ws.Range("N13").Value = -5280
# Row 132
$ws.Range("H132").Value = 5743.7
$ws.Range("I132").Value = 4626
$ws.Range("J132").Value = 6488.8335
$ws.Range("K132").Value = 13878
$ws.Range("L132").Value = 19466.5005
$ws.Range("M132").Value = -11348
$ws.Range("N132").Value = -24526.5005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 2669.5
$ws.Range("I132").Value = 1907.5652
$ws.Range("J132").Value = 3504
$ws.Range("K132").Value = 5722.6956
$ws.Range("L132").Value = 10512
$ws.Range("M132").Value = -3192.6956
$ws.Range("N132").Value = -15572
